# fix loi import anh
# The "Hang" (Brand) column for the two Apple MacBook rows was wrongly
# showing "Macbook" - correct it to "Apple" (also fixes the now-unused
# "Macbook" shared string / bad cross references introduced while the
# image import script ran).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J14").Value = "Apple"
$ws.Range("J15").Value = "Apple"

# Selection cursor moved on to the next row while reviewing the fix
$ws.Range("I19").Select()

# A default page setup got attached to the sheet as part of the re-save
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
